{"js": "// Remove the \"Auflage\" table row (label cell \"Auflage\" + its\n// \"{{ visa_requirement }}\" value cell) from the document, as described by\n// the commit \"feat: small fixes and removes auflage\".\n\n// Locate the row by searching for the unique \"Auflage\" label text rather\n// than hard-coding table/row indices, so the script is resilient to minor\n// structural differences elsewhere in the document.\nconst results = context.document.body.search(\"Auflage\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the \"Auflage\" table row to remove.');\n}\n\nconst found = results.items[0];\nconst cell = found.parentTableCell;\ncell.load(\"rowIndex\");\nconst table = cell.parentTable;\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst row = table.rows.items[cell.rowIndex];\nrow.delete();\nawait context.sync();\n", "ps1": "# Remove the \"Auflage\" table row (label cell \"Auflage\" + its\n# \"{{ visa_requirement }}\" value cell) from the document, as described by\n# the commit \"feat: small fixes and removes auflage\".\n\n$d = $word.ActiveDocument\n\n# Locate the row by its unique \"Auflage\" label text rather than a\n# hard-coded table/row index, so the script is resilient to minor\n# structural differences elsewhere in the document.\n$targetRow = $null\nforeach ($t in $d.Tables) {\n    foreach ($row in $t.Rows) {\n        if ($row.Range.Text.Contains(\"Auflage\")) {\n            $targetRow = $row\n            break\n        }\n    }\n    if ($targetRow -ne $null) { break }\n}\n\nif ($targetRow -eq $null) {\n    throw 'Could not find the \"Auflage\" table row to remove.'\n}\n\n$targetRow.Delete()\n"}
